$d = $word.ActiveDocument

# --- Edit 1: insert a new bullet paragraph after the "...60%." bullet ---
# (right before the "Projects" Heading1 paragraph)
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "*reducing quality check and content packaging time by 60%*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $p = $d.Paragraphs.Item($targetIndex)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "• Configured and Deployed File based Media Asset Management Playout systems for Linear and Non-Linear workflows."
}

# --- Edit 2: remove the trailing whitespace-only paragraph at the very end ---
# (the last paragraph in the document body, right before the final sectPr)
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
if ($lastPara.Range.Text.Trim() -eq "") {
    $lastPara.Range.Delete()
}
